$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 edits -----------------------------------------------------
# A9 becomes a formula referencing the row above (value stays 7)
$ws.Range("A9").Formula = "=A8+1"
# C9: Year_end 2017 -> 2018
$ws.Range("C9").Value = 2018

# --- New string values, written in the order the original author must
# have typed them so the shared-strings table comes out in the same
# order as the target workbook. ----------------------------------------
$ws.Range("E11").Value = "PI of Modeling Core"
$ws.Range("C10").Value = "present"
$ws.Range("D10").Value = "Adoption of innovations"
$ws.Range("F10").Value = "Leslie and John (Mac) McQuown Gift"
$ws.Range("F11").Value = "National Institute of Alergy and Infectious Diseases"
$ws.Range("D11").Value = "Successful Clinical Response in Pneumonia Therapy (SCRIPT) Systems Biology Center"

# --- Remaining row 10 / row 11 values ---------------------------------
# Column A keeps the running index via a formula, formatted like the
# rows above it (bold, bordered, centered style).
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Formula = "=A9+1"

$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Formula = "=A10+1"

$ws.Range("B10").Value = 2014
$ws.Range("E10").Value = "PI"

$ws.Range("B11").Value = 2018
$ws.Range("C11").Value = 2023

# --- Selection ---------------------------------------------------------
$ws.Range("D12").Select()
